$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Charles"
$ws.Range("B5").Value = 5
$ws.Range("G5").Value = 5

$ws.Range("A6").Value = "Matheus"
$ws.Range("B6").Value = 4.03
$ws.Range("G6").Value = 4.03

$ws.Range("A7").Value = "Murilo"
$ws.Range("B7").Value = 3.61
$ws.Range("G7").Value = 3.61

$ws.Range("A8").Value = "Teste1"
$ws.Range("B8").Value = 3.24
$ws.Range("G8").Value = 3.24

$ws.Range("A9").Value = "teste2"
$ws.Range("B9").Value = 4.38
$ws.Range("G9").Value = 4.38

$ws.Range("A10").Value = "teste"
$ws.Range("B10").Value = 5
$ws.Range("G10").Value = 5

$ws.Range("A11").Value = "teste"
$ws.Range("B11").Value = 3.96
$ws.Range("G11").Value = 3.96

$ws.Range("A12").Value = "uy"
$ws.Range("B12").Value = 1.64
$ws.Range("G12").Value = 1.64

$ws.Range("A13").Value = "teste"
$ws.Range("B13").Value = 3.59
$ws.Range("G13").Value = 3.59

$ws.Range("A14").Value = "Guilherme Ormond"
$ws.Range("B14").Value = 3.57
$ws.Range("C14").Value = ""
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = ""
$ws.Range("F14").Value = ""
$ws.Range("G14").Value = 3.57

$ws.Range("A15").Value = "teste"
$ws.Range("B15").Value = 3.67
$ws.Range("G15").Value = 3.67

$ws.Range("A16").Value = "danilo"
$ws.Range("B16").Value = 3.97
$ws.Range("G16").Value = 3.97

$ws.Range("A17").Value = "Murilo"
$ws.Range("B17").Value = 4.55
$ws.Range("C17").Value = ""
$ws.Range("D17").Value = ""
$ws.Range("E17").Value = ""
$ws.Range("F17").Value = ""
$ws.Range("G17").Value = 4.55

$ws.Range("A18").Value = "Ricardo Scopel"
$ws.Range("B18").Value = 3.73
$ws.Range("G18").Value = 3.73

$ws.Range("A19").Value = "Gustavo Bertoluzzi Cardoso"
$ws.Range("B19").Value = 3.72
$ws.Range("G19").Value = 3.72

$ws.Range("A20").Value = "Gustavo Flores"
$ws.Range("B20").Value = 3.38
$ws.Range("G20").Value = 3.38

$ws.Range("A21").Value = "Luana"
$ws.Range("B21").Value = 3.56
$ws.Range("G21").Value = 3.56

$ws.Range("A22").Value = "Rafael Testa"
$ws.Range("B22").Value = 4.49
$ws.Range("G22").Value = 4.49

$ws.Range("A23").Value = "Rafael"
$ws.Range("B23").Value = 4.48
$ws.Range("G23").Value = 4.48

$ws.Range("A24").Value = "teste4"
$ws.Range("B24").Value = 3.75
$ws.Range("G24").Value = 3.75

$ws.Range("A25").Value = "joao"
$ws.Range("B25").Value = 3.66
$ws.Range("G25").Value = 3.66

$ws.Range("A26").Value = "Gabriel"
$ws.Range("B26").Value = 4.03
$ws.Range("G26").Value = 4.03

$ws.Range("A27").Value = "Teste5"
$ws.Range("B27").Value = 3.9
$ws.Range("G27").Value = 3.9
